# Add a "Total:" label with a jxls sum formula placeholder below the
# existing data table on Sheet1, matching the style of the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New row 14: D14 = "Total:", E14 = "$[SUM(E11)]" (jxls command, stored as text)
# Shared-string table order matches authoring order: the formula placeholder
# string is interned before the "Total:" label, so write E14 first.
$ws.Range("E14").Value = "$[SUM(E11)]"
$ws.Range("D14").Value = "Total:"

# Match the bold header style used by the table header row (D9:E9).
$ws.Range("D14:E14").Font.Bold = $true

# Update the active selection to reflect the newly-added cells.
$ws.Range("D14:E14").Select()
